# Updates cryptocurrency price (D) and 1h volume % (E) columns to match the
# latest scrape. Values that look numeric (e.g. "0.9997") must be forced back
# to text so Excel does not silently convert them to a Number cell -- the
# source data stores these as plain strings (note some "prices" even contain
# multiple "." separators, e.g. "27.236.09", which are not valid numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.236.09"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "1.904.60"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9997"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "306.06"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9993"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5386"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  +3.17%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3804"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  +1.24%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07286"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "22.20"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  +5.30%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08172"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "95.79"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.338"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.9995"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("E16").Value = "  +2.27%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.000008671"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.9994"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "27.255.61"
$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").Value = "1.130.20"
$ws.Range("E20").Value = "  -40.39%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.042"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "10.81"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.519"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "148.58"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.306"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("E26").Value = "  +1.28%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.754"
$ws.Cells.Item(27, 4).Style = "Normal"

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "116.45"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "4.850"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  +1.47%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.707"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -3.88%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.09206"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  -0.09%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.8307"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  +5.26%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05076"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.998"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.326"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  -2.66%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.676"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  +3.59%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.5914"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  +4.28%  "

$ws.Range("E39").Value = "  +0.82%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.082"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "9.294"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  +3.33%  "

$ws.Range("E42").Value = "  +1.68%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "116.60"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.5115"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  +5.23%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.1530"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.24"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  +1.91%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.9989"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("E48").Value = "  +1.38%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "38.28"
$ws.Cells.Item(49, 4).Style = "Normal"

$ws.Range("E50").Value = "  +3.02%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "63.46"
$ws.Cells.Item(51, 4).Style = "Normal"
